$wb = $excel.ActiveWorkbook

# Target sheet: "CUMPLIMIENTO MENSUAL" (3rd worksheet, holds the PRESUPUESTO/VENTA/POR CUMPLIR/CUMPLIMIENTO table)
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3 (PORCELANATO group)
$ws.Range("D3").Value = 548.46
$ws.Range("E3").Value = 16951.54
$ws.Range("F3").Value = 0.03134057142857143

# Row 4 (TOTAL)
$ws.Range("D4").Value = 4368.41
$ws.Range("E4").Value = 13131.59
$ws.Range("F4").Value = 0.2496234285714286
